# Integrate Agora data for elec/BPaFF
#
# Updates the Boolean Peaking and Flexibility Flags for several plant
# types across the two data sheets, then restores the selection/active
# cell state that Excel records for each sheet.

$wb = $excel.ActiveWorkbook

# --- BPaFF-BITPTaP: "Is This Plant Type a Peaker" -------------------------
$wsPeaker = $wb.Worksheets.Item("BPaFF-BITPTaP")
# petroleum is no longer flagged as a peaker
$wsPeaker.Range("B11").Value = 0

# --- BPaFF-BDTPTPF: "Does This Plant Type Provide Flexibility" -----------
$wsFlex = $wb.Worksheets.Item("BPaFF-BDTPTPF")
# hard coal now provides flexibility
$wsFlex.Range("B2").Value = 1
# biomass now provides flexibility
$wsFlex.Range("B9").Value = 1

# --- Restore per-sheet selection / active cell ----------------------------
$wsPeaker.Activate()
$wsPeaker.Range("B17").Select()

$wsFlex.Activate()
$wsFlex.Range("B17").Select()

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("B4").Select()
